# Updates the cryptos list price/volume columns (and fixes the ARBITRUM/Hedera
# row order) to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    # The source cells are plain text (e.g. "42.876.31", "  +0.26%  ").
    # Assigning a numeric-looking string via .Value would make Excel
    # coerce it to a number, so force Text format for the write and
    # then restore the original (unstyled) cell style afterwards.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '42.853.69'
Set-TextCell "E2" '  +0.22%  '

Set-TextCell "D3" '2.567.99'
Set-TextCell "E3" '  +1.53%  '

Set-TextCell "E4" '  +0.09%  '

Set-TextCell "D5" '313.33'
Set-TextCell "E5" '  -0.60%  '

Set-TextCell "D6" '99.51'
Set-TextCell "E6" '  +3.90%  '

Set-TextCell "E7" '  -0.21%  '

Set-TextCell "E8" '  -0.01%  '

Set-TextCell "E9" '  +0.21%  '

Set-TextCell "D10" '35.73'

Set-TextCell "D11" '0.0812'
Set-TextCell "E11" '  +0.25%  '

Set-TextCell "D12" '7.47'
Set-TextCell "E12" '  -1.19%  '

Set-TextCell "D13" '2.963.60'
Set-TextCell "E13" '  +1.61%  '

Set-TextCell "E14" '  -1.07%  '

Set-TextCell "D15" '15.77'
Set-TextCell "E15" '  +3.81%  '

Set-TextCell "D16" '2.559.06'
Set-TextCell "E16" '  +0.43%  '

Set-TextCell "D17" '0.846'
Set-TextCell "E17" '  -0.84%  '

Set-TextCell "D18" '42.892.52'
Set-TextCell "E18" '  +0.11%  '

Set-TextCell "D19" '6.76'
Set-TextCell "E19" '  -1.44%  '

Set-TextCell "E20" '  -3.01%  '

Set-TextCell "D21" '0.0₃0963'
Set-TextCell "E21" '  -0.03%  '

Set-TextCell "D22" '69.57'
Set-TextCell "E22" '  -0.48%  '

Set-TextCell "E23" '  -1.43%  '

Set-TextCell "D24" '2.94'
Set-TextCell "E24" '  +0.12%  '

Set-TextCell "D25" '2.08'
Set-TextCell "E25" '  -0.14%  '

Set-TextCell "D26" '27.03'
Set-TextCell "E26" '  +1.17%  '

Set-TextCell "E27" '  -0.62%  '

Set-TextCell "E28" '  -1.27%  '

Set-TextCell "D29" '40.11'
Set-TextCell "E29" '  -0.80%  '

Set-TextCell "E30" '  -1.43%  '

Set-TextCell "D31" '156.97'
Set-TextCell "E31" '  -0.63%  '

Set-TextCell "E32" '  -1.70%  '

Set-TextCell "E33" '  +0.75%  '

Set-TextCell "B34" 'Hedera'
Set-TextCell "C34" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell "D34" '0.0800'
Set-TextCell "E34" '  +2.20%  '

Set-TextCell "B35" 'ARBITRUM'
Set-TextCell "C35" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell "D35" '2.11'
Set-TextCell "E35" '  -3.03%  '

Set-TextCell "D37" '18.71'
Set-TextCell "E37" '  -1.85%  '

Set-TextCell "D38" '2.57'
Set-TextCell "E38" '  +11.31%  '

Set-TextCell "E39" '  -0.31%  '

Set-TextCell "E40" '  -0.13%  '

Set-TextCell "D41" '23.33'
Set-TextCell "E41" '  +0.36%  '

Set-TextCell "D42" '4.14'
Set-TextCell "E42" '  +7.91%  '

Set-TextCell "E43" '  -0.09%  '

Set-TextCell "D44" '0.0302'
Set-TextCell "E44" '  -0.53%  '

Set-TextCell "D45" '3.24'
Set-TextCell "E45" '  -1.89%  '

Set-TextCell "D46" '2.006.99'
Set-TextCell "E46" '  -1.10%  '

Set-TextCell "D47" '9.01'
Set-TextCell "E47" '  -2.07%  '

Set-TextCell "D48" '2.813.44'
Set-TextCell "E48" '  +1.47%  '

Set-TextCell "D49" '0.198'
Set-TextCell "E49" '  +3.18%  '

Set-TextCell "D50" '81.80'
Set-TextCell "E50" '  -3.76%  '

Set-TextCell "D51" '74.55'
Set-TextCell "E51" '  -0.51%  '
